# "Generate Report for Handback"
#
# The handback report generator re-ran for the 048510b4-... source file
# (row 2 on each localized-language sheet / row 2 in the Overview table),
# refreshing its timestamps. The c0635734-... file's row (row 3) was not
# regenerated this time and keeps its previous timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# "Latest HO Xliff Generate Date" for 048510b4-... (row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 08:51:53"

# --- zh-cn sheet -------------------------------------------------------
# "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the 048510b4-... row (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 08:51:47"
$wsZhCn.Range("K2").Value = "2016-08-25 08:52:14"

# --- de-de sheet -------------------------------------------------------
# "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the 048510b4-... row (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 08:51:53"
$wsDeDe.Range("K2").Value = "2016-08-25 08:52:22"
